$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.217.60"
$ws.Range("E2").Value = "  -0.37%  "
$ws.Range("D3").Value = "1.826.07"
$ws.Range("E3").Value = "  -0.80%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").Value = "'236.42"
$ws.Range("E5").Value = "  -1.51%  "
$ws.Range("D6").Value = "'0.6034"
$ws.Range("E6").Value = "  -3.81%  "
$ws.Range("D7").Value = "'1.003"
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("D8").Value = "'0.07139"
$ws.Range("E8").Value = "  -3.77%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "'0.2816"
$ws.Range("E9").Value = "  -2.71%  "
$ws.Range("B10").Value = "Solana"
$ws.Range("C10").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D10").Value = "'24.11"
$ws.Range("E10").Value = "  -2.74%  "
$ws.Range("D11").Value = "'0.07659"
$ws.Range("E11").Value = "  -0.94%  "
$ws.Range("D12").Value = "1.859.25"
$ws.Range("E12").Value = "  +1.20%  "
$ws.Range("D13").Value = "'4.798"
$ws.Range("E13").Value = "  -3.62%  "
$ws.Range("D14").Value = "'0.6433"
$ws.Range("E14").Value = "  -4.94%  "
$ws.Range("D15").Value = "'0.000009808"
$ws.Range("E15").Value = "  -3.85%  "
$ws.Range("D16").Value = "'79.62"
$ws.Range("E16").Value = "  -2.84%  "
$ws.Range("D17").Value = "2.038.37"
$ws.Range("E17").Value = "  -2.37%  "
$ws.Range("D18").Value = "'6.031"
$ws.Range("E18").Value = "  -3.36%  "
$ws.Range("D19").Value = "29.204.58"
$ws.Range("E19").Value = "  -0.45%  "
$ws.Range("D20").Value = "'231.39"
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("D21").Value = "'1.002"
$ws.Range("E21").Value = "  +0.24%  "
$ws.Range("D22").Value = "'11.80"
$ws.Range("E22").Value = "  -3.95%  "
$ws.Range("D23").Value = "'7.038"
$ws.Range("E23").Value = "  -5.06%  "
$ws.Range("D24").Value = "'1.003"
$ws.Range("E24").Value = "  +0.30%  "
$ws.Range("D25").Value = "'156.03"
$ws.Range("E25").Value = "  -1.81%  "
$ws.Range("D26").Value = "'8.099"
$ws.Range("E26").Value = "  -4.39%  "
$ws.Range("D27").Value = "'0.1284"
$ws.Range("E27").Value = "  -5.02%  "
$ws.Range("D28").Value = "'16.75"
$ws.Range("E28").Value = "  -3.83%  "
$ws.Range("D29").Value = "'0.06814"
$ws.Range("E29").Value = "  +4.02%  "
$ws.Range("D30").Value = "'1.466"
$ws.Range("E30").Value = "  +1.72%  "
$ws.Range("D31").Value = "'1.460"
$ws.Range("E31").Value = "  -1.38%  "
$ws.Range("D32").Value = "'3.835"
$ws.Range("E32").Value = "  -5.69%  "
$ws.Range("D33").Value = "'3.773"
$ws.Range("E33").Value = "  -7.14%  "
$ws.Range("D34").Value = "'1.134"
$ws.Range("E34").Value = "  -0.34%  "
$ws.Range("D35").Value = "'1.715"
$ws.Range("E35").Value = "  -6.38%  "
$ws.Range("D36").Value = "'0.6601"
$ws.Range("E36").Value = "  -4.66%  "
$ws.Range("D37").Value = "'2.537"
$ws.Range("E37").Value = "  -1.12%  "
$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "1.225.20"
$ws.Range("E38").Value = "  -1.49%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").Value = "'2.762"
$ws.Range("E39").Value = "  -2.15%  "
$ws.Range("E40").Value = "  -4.49%  "
$ws.Range("D41").Value = "'6.567"
$ws.Range("E41").Value = "  -2.58%  "
$ws.Range("D42").Value = "'0.9265"
$ws.Range("E42").Value = "  -0.68%  "
$ws.Range("D43").Value = "'1.002"
$ws.Range("E43").Value = "  +0.28%  "
$ws.Range("D44").Value = "1.961.80"
$ws.Range("E44").Value = "  -1.45%  "
$ws.Range("D45").Value = "'99.69"
$ws.Range("E45").Value = "  -0.99%  "
$ws.Range("D46").Value = "'63.41"
$ws.Range("E46").Value = "  -3.39%  "
$ws.Range("E47").Value = "  +1.10%  "
$ws.Range("D48").Value = "'1.634"
$ws.Range("E48").Value = "  -4.60%  "
$ws.Range("D49").Value = "'6.598"
$ws.Range("E49").Value = "  -6.28%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D50").Value = "'8.534"
$ws.Range("E50").Value = "  -5.12%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").Value = "'0.05593"
$ws.Range("E51").Value = "  -1.45%  "
